# Moved speed function to new file, added motor calibration
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Populate Sheet2 with the f -> v lookup table --------------------------
$ws2.Range("A1").Value = "f"
$ws2.Range("B1").Value = "v"

$row = 2
for ($x = 0; $x -le 2800; $x += 100) {
    $ws2.Cells.Item($row, 1).Value = $x
    $formula = "=(A$row-Sheet1!T`$32)/Sheet1!T`$31"
    $ws2.Cells.Item($row, 2).Formula = $formula
    $row++
}

# --- Sheet1 view / selection state -----------------------------------------
$null = $ws1.Activate()
$null = $ws1.Range("C31").Select()
try { $excel.ActiveWindow.ScrollColumn = 2 } catch {}

# --- Page setup on Sheet1 (printed portrait) --------------------------------
$ws1.PageSetup.Orientation = 1

# --- Sheet2 becomes the active sheet, header row selected -------------------
$null = $ws2.Activate()
$null = $ws2.Range("A1:B1").Select()
